$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on column D (prices) so Excel does not
# reinterpret the values (e.g. "1.00", "0.0346", trailing zeros,
# "." used as thousands separator) as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.365.97'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '3.117.09'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '576.43'
$ws.Range('D6').Value = '179.08'
$ws.Range('E6').Value = '  +6.57%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.116.95'
$ws.Range('E8').Value = '  +1.44%  '
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('E10').Value = '  +2.52%  '
$ws.Range('D12').Value = '0.468'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '0.0000242'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').Value = '36.61'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('D16').Value = '3.632.44'
$ws.Range('E16').Value = '  +1.25%  '
$ws.Range('D17').Value = '67.341.93'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').Value = '3.112.97'
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('D20').Value = '16.47'
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').Value = '485.19'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').Value = '0.689'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '7.72'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').Value = '83.65'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = '2.29'
$ws.Range('E25').Value = '  +3.80%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '12.73'
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').Value = '10.55'
$ws.Range('E27').Value = '  +2.76%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '8.00'
$ws.Range('E29').Value = '  +2.27%  '
$ws.Range('D30').Value = '2.33'
$ws.Range('E30').Value = '  +2.23%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').Value = '  +1.71%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').Value = '0.0₃0941'
$ws.Range('E34').Value = '  +3.46%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '48.39'
$ws.Range('E36').Value = '  +4.59%  '
$ws.Range('D37').Value = '0.953'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = '5.60'
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('D39').Value = '0.321'
$ws.Range('E39').Value = '  +6.19%  '
$ws.Range('D40').Value = '49.25'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').Value = '8.32'
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').Value = '2.70'
$ws.Range('E44').Value = '  +8.31%  '
$ws.Range('D45').Value = '2.792.73'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('D46').Value = '374.12'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '26.76'
$ws.Range('E47').Value = '  +9.56%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0346'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = '135.69'
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  +10.22%  '
